# Generate Report for Handoff
# Replaces the previously-handed-back file pair (365764c5...md / a997d19f...md)
# with a freshly handed-off file pair (0a48022b...md / ffff00348e7b...md),
# clearing the stale "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" data and stamping new handoff timestamps + xliff names.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "365764c5-d128-40bc-9cee-edb6cb33f643"
$oldUuid2 = "a997d19f-6a67-4018-8d32-d9177a7f1463"
$newUuid1 = "0a48022b-3f22-4b64-95fa-057cae1d5fe7"
$newUuid2 = "ffff00348e7b-4708-4163-971f-9b27c43aeefa"

$newStatus = "Ready for handoff"
$newHoDate = "2016-08-31 21:18:23"
$newHandoffDateTime = "2016-08-31 21:18:18"
$clearedHandbackDateTime = "0001-01-01 00:00:00"

$newXlfHashZhCn = "1eae068af0547b479fed30e3d824c7431b952718"
$newXlfHashDeDe = "1eae068af0547b479fed30e3d824c7431b952718"

$newHandoffFileZhCn = "$newUuid1.$newXlfHashZhCn.zh-cn.xlf"
$newHandoffFileDeDe = "$newUuid1.$newXlfHashDeDe.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$targetB2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff27aa36a9b5cac05d4139de009495f5851d87ec/e2e/$oldUuid1.md"
$targetB3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff27aa36a9b5cac05d4139de009495f5851d87ec/e2e/$oldUuid2.md"

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $targetB2, "", "", "e2e\$newUuid1.md")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $targetB3, "", "", "e2e\$newUuid2.md")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newHoDate

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$targetA2ZhCn = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff27aa36a9b5cac05d4139de009495f5851d87ec/e2e/$oldUuid1.md"
$targetA3ZhCn = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff27aa36a9b5cac05d4139de009495f5851d87ec/e2e/$oldUuid2.md"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $targetA2ZhCn, "", "", "$newUuid1.md")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("G2").Value = $newHandoffFileZhCn
$wsZhCn.Range("H2").Value = $newHandoffDateTime
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("I2").Hyperlinks.Delete()
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = $clearedHandbackDateTime

$wsZhCn.Range("A3").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $targetA3ZhCn, "", "", "$newUuid2.md")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = $newHandoffFileZhCn
$wsZhCn.Range("H3").Value = $newHandoffDateTime
$wsZhCn.Range("I3").Style = "Normal"
$wsZhCn.Range("I3").Hyperlinks.Delete()
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = $clearedHandbackDateTime

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$targetA2DeDe = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff27aa36a9b5cac05d4139de009495f5851d87ec/e2e/$oldUuid1.md"
$targetA3DeDe = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff27aa36a9b5cac05d4139de009495f5851d87ec/e2e/$oldUuid2.md"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $targetA2DeDe, "", "", "$newUuid1.md")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("G2").Value = $newHandoffFileDeDe
$wsDeDe.Range("H2").Value = $newHoDate
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("I2").Hyperlinks.Delete()
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = $clearedHandbackDateTime

$wsDeDe.Range("A3").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $targetA3DeDe, "", "", "$newUuid2.md")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = $newHandoffFileDeDe
$wsDeDe.Range("H3").Value = $newHoDate
$wsDeDe.Range("I3").Style = "Normal"
$wsDeDe.Range("I3").Hyperlinks.Delete()
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = $clearedHandbackDateTime
